$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$data = New-Object "object[,]" 30,8

$data[0,0] = 0
$data[0,1] = "walkingToRunning"
$data[0,2] = -4.399603349646254
$data[0,3] = -32.20032115818277
$data[0,4] = -1.038676389713888
$data[0,5] = 2.395846366882324
$data[0,6] = 5.165933609008789
$data[0,7] = -3.433014392852783

$data[1,0] = 100
$data[1,1] = "walkingToRunning"
$data[1,2] = 2.827775232570724
$data[1,3] = -30.86045528687148
$data[1,4] = 2.059290839839234
$data[1,5] = -4.131561756134033
$data[1,6] = 1.293697118759155
$data[1,7] = 3.305315971374512

$data[2,0] = 200
$data[2,1] = "walkingToRunning"
$data[2,2] = 9.072769558306828
$data[2,3] = -14.51345410297827
$data[2,4] = 10.26892458532274
$data[2,5] = -11.94900512695312
$data[2,6] = -0.5390903949737549
$data[2,7] = 0.9473530054092408

$data[3,0] = 300
$data[3,1] = "walkingToRunning"
$data[3,2] = 1.526257913137283
$data[3,3] = -15.38509752332542
$data[3,4] = -1.262450129715198
$data[3,5] = 6.543986320495605
$data[3,6] = -9.976800918579102
$data[3,7] = -2.196774959564209

$data[4,0] = 400
$data[4,1] = "walkingToRunning"
$data[4,2] = 6.41234910119443
$data[4,3] = -10.75634378256254
$data[4,4] = -6.548929008011939
$data[4,5] = 2.871486902236938
$data[4,6] = -5.876065254211426
$data[4,7] = -3.124087810516357

$data[5,0] = 500
$data[5,1] = "walkingToRunning"
$data[5,2] = 32.94686977150512
$data[5,3] = -15.4527910517663
$data[5,4] = -8.939652902563815
$data[5,5] = 7.37382698059082
$data[5,6] = 11.17136192321777
$data[5,7] = -1.706753373146057

$data[6,0] = 600
$data[6,1] = "walkingToRunning"
$data[6,2] = 53.20190917339517
$data[6,3] = -41.65597895986004
$data[6,4] = -19.47473356895838
$data[6,5] = -4.830907821655273
$data[6,6] = 0.3003380000591278
$data[6,7] = 1.922602653503418

$data[7,0] = 700
$data[7,1] = "walkingToRunning"
$data[7,2] = 26.88092610516489
$data[7,3] = -29.92789514040201
$data[7,4] = -3.621464901363669
$data[7,5] = -6.295113563537598
$data[7,6] = 2.037784337997437
$data[7,7] = 3.717573165893554

$data[8,0] = 800
$data[8,1] = "walkingToRunning"
$data[8,2] = -14.28073297579263
$data[8,3] = -29.62231750586597
$data[8,4] = 40.7812210653246
$data[8,5] = -1.753358721733093
$data[8,6] = 12.93171119689941
$data[8,7] = 2.580935955047607

$data[9,0] = 900
$data[9,1] = "walkingToRunning"
$data[9,2] = 29.81173590040617
$data[9,3] = -15.35024061891273
$data[9,4] = 27.41680280449521
$data[9,5] = -3.18294358253479
$data[9,6] = -2.735665798187256
$data[9,7] = 0.6512094736099243

$data[10,0] = 1000
$data[10,1] = "walkingToRunning"
$data[10,2] = 38.08487456606819
$data[10,3] = 8.25389982253002
$data[10,4] = -6.631011874405364
$data[10,5] = 7.503256797790527
$data[10,6] = 0.8255133628845215
$data[10,7] = -1.698231339454651

$data[11,0] = 1100
$data[11,1] = "walkingToRunning"
$data[11,2] = 15.41566228866577
$data[11,3] = -42.39830100536346
$data[11,4] = 41.37434816360474
$data[11,5] = -1.000749349594116
$data[11,6] = -2.794787883758545
$data[11,7] = -1.550159573554993

$data[12,0] = 1200
$data[12,1] = "walkingToRunning"
$data[12,2] = 32.54000410099627
$data[12,3] = -39.22240943515485
$data[12,4] = 45.66644129802327
$data[12,5] = -2.628473520278931
$data[12,6] = -2.876280546188354
$data[12,7] = -0.96160089969635

$data[13,0] = 1300
$data[13,1] = "walkingToRunning"
$data[13,2] = 24.84834649390785
$data[13,3] = -8.03618379229086
$data[13,4] = 24.45961570739691
$data[13,5] = -5.674064636230469
$data[13,6] = -1.664009690284729
$data[13,7] = -3.089466571807861

$data[14,0] = 1400
$data[14,1] = "walkingToRunning"
$data[14,2] = -12.76972269274368
$data[14,3] = -22.28177534673628
$data[14,4] = 15.843390720408
$data[14,5] = 0.2152500003576278
$data[14,6] = 15.8856897354126
$data[14,7] = -3.721701145172119

$data[15,0] = 1500
$data[15,1] = "walkingToRunning"
$data[15,2] = -17.05221034571044
$data[15,3] = -9.689534740349583
$data[15,4] = -7.262813420639774
$data[15,5] = 1.613942265510559
$data[15,6] = -4.874184131622314
$data[15,7] = -0.0044607948511838

$data[16,0] = 1600
$data[16,1] = "walkingToRunning"
$data[16,2] = 7.601734495654355
$data[16,3] = -7.369463099647012
$data[16,4] = 0.05982191538072679
$data[16,5] = 4.66752290725708
$data[16,6] = -2.795853137969971
$data[16,7] = -8.249074935913086

$data[17,0] = 1700
$data[17,1] = "walkingToRunning"
$data[17,2] = -5.084137159524504
$data[17,3] = -17.42479633056011
$data[17,4] = 13.42315847357525
$data[17,5] = 1.00194776058197
$data[17,6] = -2.820886850357056
$data[17,7] = 4.437692165374756

$data[18,0] = 1800
$data[18,1] = "walkingToRunning"
$data[18,2] = -2.988111949458581
$data[18,3] = -1.541268535496201
$data[18,4] = 24.50189582588753
$data[18,5] = -8.889565467834473
$data[18,6] = 1.164800047874451
$data[18,7] = 11.45938301086426

$data[19,0] = 1900
$data[19,1] = "walkingToRunning"
$data[19,2] = 16.88818173740272
$data[19,3] = -21.43393101642941
$data[19,4] = 26.25781544950798
$data[19,5] = 2.723947763442993
$data[19,6] = 6.857040882110596
$data[19,7] = -5.117730140686035

$data[20,0] = 2000
$data[20,1] = "walkingToRunning"
$data[20,2] = 11.71483505878375
$data[20,3] = -29.24488582807863
$data[20,4] = 14.33557478914648
$data[20,5] = 1.208609104156494
$data[20,6] = -8.24947452545166
$data[20,7] = -0.5370930433273315

$data[21,0] = 2100
$data[21,1] = "walkingToRunning"
$data[21,2] = -17.89474326064904
$data[21,3] = 1.058929585918881
$data[21,4] = -6.072563363104724
$data[21,5] = 8.961604118347168
$data[21,6] = 6.756906032562256
$data[21,7] = -2.651110410690308

$data[22,0] = 2200
$data[22,1] = "walkingToRunning"
$data[22,2] = -6.588068613071894
$data[22,3] = -42.04801777711909
$data[22,4] = 27.395873059932
$data[22,5] = -2.561361789703369
$data[22,6] = 1.376787781715393
$data[22,7] = 1.834718346595764

$data[23,0] = 2300
$data[23,1] = "walkingToRunning"
$data[23,2] = 4.73087814911125
$data[23,3] = -46.72824330182412
$data[23,4] = 40.59463504909236
$data[23,5] = -2.701976776123047
$data[23,6] = 1.867341995239258
$data[23,7] = 1.237104892730713

$data[24,0] = 2400
$data[24,1] = "walkingToRunning"
$data[24,2] = 13.50610577691452
$data[24,3] = -5.020019875359202
$data[24,4] = 23.69000314928804
$data[24,5] = -5.291634559631348
$data[24,6] = -3.584681510925293
$data[24,7] = -2.499310255050659

$data[25,0] = 2500
$data[25,1] = "walkingToRunning"
$data[25,2] = -13.47701175925547
$data[25,3] = -14.91354519067346
$data[25,4] = 16.03211706692484
$data[25,5] = -0.4974119365215301
$data[25,6] = 12.873122215271
$data[25,7] = 5.119993686676025

$data[26,0] = 2600
$data[26,1] = "walkingToRunning"
$data[26,2] = -16.10707013631821
$data[26,3] = -7.06407377891983
$data[26,4] = -12.22038123533914
$data[26,5] = -0.0159123875200748
$data[26,6] = 1.14615797996521
$data[26,7] = 3.750063896179199

$data[27,0] = 2700
$data[27,1] = "walkingToRunning"
$data[27,2] = 17.60490359473454
$data[27,3] = -9.783274660405542
$data[27,4] = -4.873167308335432
$data[27,5] = 8.204200744628906
$data[27,6] = 6.466621398925781
$data[27,7] = -1.093693733215332

$data[28,0] = 2800
$data[28,1] = "walkingToRunning"
$data[28,2] = -13.56971339589527
$data[28,3] = -42.06590333918953
$data[28,4] = 1.837090728209249
$data[28,5] = -1.878527283668518
$data[28,6] = -2.402770519256592
$data[28,7] = 0.963331937789917

$data[29,0] = 2900
$data[29,1] = "walkingToRunning"
$data[29,2] = -0.4818755535735217
$data[29,3] = -32.3410521585916
$data[29,4] = 8.205087691238255
$data[29,5] = -2.405300617218018
$data[29,6] = 1.181844353675843
$data[29,7] = 4.602807998657227

$rng = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(31, 8))
$rng.Value = $data